$wb = $excel.ActiveWorkbook

# --- hubs sheet: add the new "adjustment time for total utilized capacity" row,
#     a fix so the capacity adjusts (fills) smoothly instead of instantly ---
$hubs = $wb.Worksheets.Item("hubs")
$hubs.Activate()
$hubs.Range("A4").Value = "adjustment time for total utilized capacity"
$hubs.Range("B4").Value = 24
$hubs.Range("C4").Value = "Month"
$hubs.Range("A5").Select() | Out-Null

# --- leave the "demand" sheet as the active / selected sheet & cell ---
$demand = $wb.Worksheets.Item("demand")
$demand.Activate()
$demand.Range("B1").Select() | Out-Null
